$wb = $excel.ActiveWorkbook

# GLOBAL RESULTS
$ws = $wb.Worksheets.Item("GLOBAL RESULTS")
$ws.Range("C6").Value = 62876.30628779709
$ws.Range("C7").Value = 56588.675659017375
$ws.Range("C9").Value = 11633.190117226004
$ws.Range("C10").Value = 54778.116170571084
$ws.Range("C11").Value = 51243.116170571084
$ws.Range("C12").Value = 17035.0
$ws.Range("C13").Value = 13500.0
$ws.Range("C14").Value = 37743.11617057109
$ws.Range("C15").Value = 36969.622526539715
$ws.Range("C16").Value = 314.40635303136366
$ws.Range("C18").Value = 2020.9499999999998
$ws.Range("C19").Value = 35263.0788795711
$ws.Range("C20").Value = 20602.655422083808

# FUSELAGE
$ws = $wb.Worksheets.Item("FUSELAGE")
$ws.Range("C2").Value = 6288.127060627273
$ws.Range("C3").Value = 8132.166666666666
$ws.Range("D3").Value = 29.325737032028755
$ws.Range("C5").Value = 8132.166666666666
$ws.Range("C8").Value = 14190.0
$ws.Range("D8").Value = 125.66337898688195
$ws.Range("C9").Value = 10196.0
$ws.Range("D9").Value = 62.14685075054604
$ws.Range("C10").Value = 6416.0
$ws.Range("D10").Value = 2.0335616335330875
$ws.Range("C11").Value = 21031.0
$ws.Range("D11").Value = 234.45570989944426
$ws.Range("C12").Value = 10802.0
$ws.Range("D12").Value = 71.78406059311477
$ws.Range("D13").Value = 1.7155019027552412
$ws.Range("C14").Value = 7092.0
$ws.Range("D14").Value = 12.783980533824291
$ws.Range("C15").Value = 7891.0
$ws.Range("D15").Value = 25.49046677839925

# WING
$ws = $wb.Worksheets.Item("WING")
$ws.Range("C2").Value = 6665.414684264909
$ws.Range("C3").Value = 7112.75
$ws.Range("D3").Value = 6.711290098590839
$ws.Range("C5").Value = 7112.749999999999
$ws.Range("A8").Value = "TORENBEEK_2013"
$ws.Range("C8").Value = 8327.0
$ws.Range("D8").Value = 24.928461235241798
$ws.Range("A9").Value = "TORENBEEK_1982"
$ws.Range("C9").Value = 6559.0
$ws.Range("D9").Value = -1.5965200862314213
$ws.Range("A10").Value = "KROO"
$ws.Range("C10").Value = 7483.0
$ws.Range("D10").Value = 12.266083274086032
$ws.Range("A11").Value = "JENKINSON"
$ws.Range("C11").Value = 6082.0
$ws.Range("D11").Value = -8.752864028732963

# HORIZONTAL TAIL
$ws = $wb.Worksheets.Item("HORIZONTAL TAIL")
$ws.Range("C2").Value = 723.1346119721363
$ws.Range("C3").Value = 795.375
$ws.Range("D3").Value = 9.989894942360628
$ws.Range("C5").Value = 795.3749999999999
$ws.Range("A8").Value = "HOWE"
$ws.Range("C8").Value = 394.0
$ws.Range("D8").Value = -45.51498524936025
$ws.Range("C9").Value = 1523.0
$ws.Range("D9").Value = 110.61085651072167
$ws.Range("A10").Value = "TORENBEEK_2013"
$ws.Range("C10").Value = 502.0
$ws.Range("D10").Value = -30.58000658674834
$ws.Range("A11").Value = "NICOLAI_1984"
$ws.Range("C11").Value = 700.0
$ws.Range("D11").Value = -3.199212371959834
$ws.Range("A12").Value = "ROSKAM"
$ws.Range("C12").Value = 52.0
$ws.Range("D12").Value = -92.8090843476313
$ws.Range("A13").Value = "NICOLAI_2013"
$ws.Range("C13").Value = 1415.0
$ws.Range("D13").Value = 95.67587784810976
$ws.Range("A14").Value = "SADRAEY"
$ws.Range("C14").Value = 1040.0
$ws.Range("D14").Value = 43.81831304737396
$ws.Range("A15").Value = "KROO"
$ws.Range("C15").Value = 737.0
$ws.Range("D15").Value = 1.9174006883794317

# VERTICAL TAIL
$ws = $wb.Worksheets.Item("VERTICAL TAIL")
$ws.Range("C2").Value = 723.1346119721363
$ws.Range("C3").Value = 672.4285714285713
$ws.Range("D3").Value = -7.011978088737747
$ws.Range("C5").Value = 672.4285714285713
$ws.Range("A8").Value = "RAYMER"
$ws.Range("C8").Value = 1523.0
$ws.Range("D8").Value = 110.61085651072167
$ws.Range("A9").Value = "TORENBEEK_2013"
$ws.Range("C9").Value = 179.0
$ws.Range("D9").Value = -75.24665573511544
$ws.Range("A10").Value = "NICOLAI_1984"
$ws.Range("C10").Value = 502.0
$ws.Range("D10").Value = -30.58000658674834
$ws.Range("C11").Value = 124.0
$ws.Range("D11").Value = -82.85243190589003
$ws.Range("A12").Value = "NICOLAI_2013"
$ws.Range("C12").Value = 1145.0
$ws.Range("D12").Value = 58.338431191579986
$ws.Range("A13").Value = "SADRAEY"
$ws.Range("C13").Value = 749.0
$ws.Range("D13").Value = 3.5768427620029772
$ws.Range("A14").Value = "KROO"
$ws.Range("C14").Value = 485.0
$ws.Range("D14").Value = -32.93088285771503

# NACELLES
$ws = $wb.Worksheets.Item("NACELLES")
$ws.Range("C2").Value = 1194.7441415191815
$ws.Range("D3").Value = 16.398701474138132
$ws.Range("A10").Value = "RAYMER"
$ws.Range("D10").Value = 15.003702654937063
$ws.Range("D11").Value = 16.17550166306597
$ws.Range("A12").Value = "NICOLAI_1984"
$ws.Range("D12").Value = 18.016900104411395
$ws.Range("A17").Value = "RAYMER"
$ws.Range("D17").Value = 15.003702654937063
$ws.Range("D18").Value = 16.17550166306597
$ws.Range("A19").Value = "NICOLAI_1984"
$ws.Range("D19").Value = 18.016900104411395

# POWER PLANT
$ws = $wb.Worksheets.Item("POWER PLANT")
$ws.Range("C2").Value = 5219.145460320637
$ws.Range("D3").Value = 23.60899657579186
$ws.Range("A11").Value = "KUNDU"
$ws.Range("C11").Value = 3265.0
$ws.Range("D11").Value = 25.116267589116617
$ws.Range("A12").Value = "ROSKAM"
$ws.Range("C12").Value = 2954.0
$ws.Range("D12").Value = 13.19860779732021
$ws.Range("A13").Value = "JENKINSON"
$ws.Range("C13").Value = 3458.0
$ws.Range("D13").Value = 32.51211434093882
$ws.Range("A18").Value = "KUNDU"
$ws.Range("C18").Value = 3265.0
$ws.Range("D18").Value = 25.116267589116617
$ws.Range("A19").Value = "ROSKAM"
$ws.Range("C19").Value = 2954.0
$ws.Range("D19").Value = 13.19860779732021
$ws.Range("A20").Value = "JENKINSON"
$ws.Range("C20").Value = 3458.0
$ws.Range("D20").Value = 32.51211434093882

# LANDING GEARS
$ws = $wb.Worksheets.Item("LANDING GEARS")
$ws.Range("C2").Value = 2578.132094857182
$ws.Range("C3").Value = 2499.2685173219097
$ws.Range("D3").Value = -3.0589424681764017
$ws.Range("C5").Value = 2499.2685173219093
$ws.Range("A9").Value = "ROSKAM"
$ws.Range("C9").Value = 2499.2685173219097
$ws.Range("D9").Value = -3.0589424681764026
$ws.Range("A11").Value = "ROSKAM"
$ws.Range("C11").Value = 386.729549170154
$ws.Range("A13").Value = "ROSKAM"
$ws.Range("C13").Value = 2112.5389681517563

# SYSTEMS
$ws = $wb.Worksheets.Item("SYSTEMS")
$ws.Range("C2").Value = 8551.852802453091
$ws.Range("C3").Value = 8209.090124153976
$ws.Range("D3").Value = -4.00805166104818
$ws.Range("C4").Value = 8209.090124153974
$ws.Range("C8").Value = 8209.090124153976
$ws.Range("D8").Value = -4.00805166104815
$ws.Range("C11").Value = 329.5722550152342
$ws.Range("C13").Value = 329.57225501523413
$ws.Range("C21").Value = 1030.8432141307385
$ws.Range("C23").Value = 1030.8432141307383
$ws.Range("C26").Value = 528.9794736539732
$ws.Range("C28").Value = 528.9794736539731
$ws.Range("C36").Value = 777.495132217799
$ws.Range("C38").Value = 777.4951322177989
$ws.Range("C41").Value = 3293.1360915880978
$ws.Range("C43").Value = 3293.1360915880973
